$wb = $excel.ActiveWorkbook

# --- Sheet "grilla de pruebas": update Capital (B3) and the divisor (B9) ---
$wsGrilla = $wb.Worksheets.Item("grilla de pruebas")

# Capital changes from 650.09 to 639 (kept as a text-like value, same as original)
$wsGrilla.Range("B3").Value = "639"

# The divisor cell B9 becomes the shared text "6" (instead of the numeric 6)
$wsGrilla.Range("B9").Value = "6"

# --- Sheet "CALCULADORA": update size/price inputs ---
$wsCalc = $wb.Worksheets.Item("CALCULADORA")

$wsCalc.Range("E3").Value = 491
$wsCalc.Range("F3").Value = 6.18
$wsCalc.Range("E4").Value = 0
$wsCalc.Range("F4").Value = 6.368
$wsCalc.Range("J6").Value = 6.093

# Recalculate so cached formula results reflect the new inputs
$wb.Application.Calculate()

# --- Update cursor/selection position on each sheet ---
$wsGrilla.Activate()
$wsGrilla.Range("A24").Select()

$wsCalc.Activate()
$wsCalc.Range("E5").Select()
